$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.675.96'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '1.834.82'
$ws.Range("E3").Value = '  +1.73%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '317.86'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5355'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.41%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3967'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +4.84%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07735'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.28%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.123'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.14%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '41.88'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.05%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '6.377'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.42%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '21.23'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.05%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.591'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.95%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '1.829.59'
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '91.94'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.00001089'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.19%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06584'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.85'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.68%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.101'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.19%  '
$ws.Range("D23").Value = '28.644.56'
$ws.Range("E23").Value = '  +1.60%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.22'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.58%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.242'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +7.43%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '20.82'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.445'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +5.35%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '157.04'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").Value = '2.042.42'
$ws.Range("E29").Value = '  +1.41%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '125.58'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +2.64%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.142'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.10%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.1119'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.91%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.767'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.99%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.662'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.61%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.07361'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.28%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.2268'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.60%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02359'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.234'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.22%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.887'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.40%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '11.47'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.86%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.6316'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.15%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.197'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.94%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.12%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.392'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.51%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '13.51'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5926'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.715'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.89%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '125.94'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.31%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.005'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.83%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.202'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.98%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06958'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.99%  '

Write-Output "Updated cryptos list"
